$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (current "Appt Time" column),
# shifting Appt Time and everything after it one column to the right.
$ws.Columns("E").Insert()

# Set the header for the newly inserted column E
$ws.Range("E1").Value = " Appt Date"

# The new column inherits the width of the column to its left (D)
$ws.Columns("E").ColumnWidth = $ws.Columns("D").ColumnWidth

# Move selection to E2 (matches the diff's new selection state)
$ws.Range("E2").Select()
